$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column X (codelist), shifting
# everything from X onward two columns to the right.
[void]$ws.Range("X1:Y1").EntireColumn.Insert()

# Give the new columns their headers.
$ws.Range("X1").Value2 = "derived_variable"
$ws.Range("Y1").Value2 = "derivation_description"

# Re-apply the AutoFilter so its range covers the now-wider table.
$ws.AutoFilterMode = $False
[void]$ws.Range("A1:AL60").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Collection_PR!_FilterDatabase") {
        $n.RefersTo = "=Collection_PR!`$A`$1:`$AL`$60"
    }
}

# Restore the active selection to match the post-edit view.
[void]$ws.Range("X3").Select()
